# Generate Report for Handoff
# Update the Priority and Latest Handoff Datetime columns for the
# "Ready for handoff" rows (rows 4-7) on both the zh-cn and de-de sheets,
# reflecting a fresh handoff xliff generation pass.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# zh-cn: rows 4-7, column E = Priority, column H = Latest Handoff Datetime
foreach ($r in 4..7) {
    $wsZh.Cells.Item($r, 5).Value = "ht"
    $wsZh.Cells.Item($r, 8).Value = "2016-08-28 10:36:47"
}

# de-de: rows 4-7, column E = Priority, column H = Latest Handoff Datetime
foreach ($r in 4..7) {
    $wsDe.Cells.Item($r, 5).Value = "ht"
    $wsDe.Cells.Item($r, 8).Value = "2016-08-28 10:36:51"
}

# Overview: rows 4-7, column G = Latest HO Xliff Generate Date
# (mirrors the de-de handoff datetime that was just refreshed)
foreach ($r in 4..7) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-08-28 10:36:51"
}
